$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("No data available") is removed; all subsequent data rows shift
# up by one, which also shrinks the used range from F94 to F93.
$ws.Rows.Item(2).Delete()
